# Update the "want to go" count (column F) on several rows across three
# worksheets to reflect newly generated data (gh-pages output).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    5  = 15511
    7  = 6
    8  = 700
    9  = 15387
    10 = 50
    11 = 8981
    12 = 373
    14 = 1010
    15 = 86
    16 = 196
    18 = 196
    20 = 46
    21 = 546
    22 = 25
    25 = 1107
    28 = 79
    30 = 39
    32 = 58
    34 = 247
    35 = 312
    36 = 447
    37 = 115
    38 = 5504
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 67

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    5  = 15511
    7  = 6
    8  = 700
    9  = 15387
    10 = 50
    11 = 8981
    12 = 373
    14 = 1010
    15 = 86
    16 = 196
    18 = 196
    20 = 46
    21 = 546
    22 = 25
    25 = 1107
    28 = 79
    30 = 39
    31 = 67
    34 = 58
    36 = 247
    37 = 312
    38 = 447
    39 = 115
    40 = 5504
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
